$d = $word.ActiveDocument

# 1. Add a third author (Gustavo Bolfe) to the author line.
$d.Content.Find.Execute("Guilherme Bolfe, Silvio Bolfe", $false, $false, $false, $false, $false, $true, 1, $false, "Guilherme Bolfe, Silvio Bolfe, Gustavo Bolfe", 2) | Out-Null

# 2. The superscript affiliation-marker line gains an extra ", " for the new author.
#    Scope the Find to that specific paragraph (3rd paragraph) so we don't touch any
#    other ", " substring in the document.
$affiliationPara = $d.Paragraphs.Item(3).Range
$affiliationPara.Find.Execute(", ", $true, $false, $false, $false, $false, $true, 1, $false, ", , ", 2) | Out-Null

# 3. Add the new author's e-mail address to the contact line.
$d.Content.Find.Execute("guilhermebolfe11@gmail.com, silviobolfe19@gmail.com", $false, $false, $false, $false, $false, $true, 1, $false, "guilhermebolfe11@gmail.com, silviobolfe19@gmail.com, bolfeguilherme@gmail.com", 2) | Out-Null

# 4. Update the "1.3. Domains" body text.
$d.Content.Find.Execute("dfasfasd.", $false, $false, $false, $false, $false, $true, 1, $false, "dfgfg.", 2) | Out-Null

# 5. Fill in the Synonyms cell (Table 3, row 2, col 2) that was empty.
$d.Tables.Item(3).Cell(2, 2).Range.Text = "ggbcvb OR gfhfgh"

# 6. Expand the search-string "(zxc)" occurrences (Table 4, both data rows) to include
#    the new synonyms.
$d.Content.Find.Execute("(zxc)", $false, $false, $false, $false, $false, $true, 1, $false, "(zxc OR ggbcvb OR gfhfgh)", 2) | Out-Null

# 7. Switch the Inclusion Rule from "Any" to "At Least".
$d.Content.Find.Execute("Inclusion Rule: Any.", $false, $false, $false, $false, $false, $true, 1, $false, "Inclusion Rule: At Least.", 2) | Out-Null

# 8. Add a new row to the last table ("Extraction Questions", Table 9) describing the
#    new "vcb" field (a Multiple Choice List with three options).
$extractionTable = $d.Tables.Item(9)
$newRow = $extractionTable.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "vcb"
$newRow.Cells.Item(2).Range.Text = "cbfghgfhfghfghfghfghfghfh"
$newRow.Cells.Item(3).Range.Text = "Multiple Choice List"
$newRow.Cells.Item(4).Range.Text = "fghfghfghfghfghfgh;" + [char]10 + "fghfghfg;" + [char]10 + "hfgh;" + [char]10
